$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (header "Förändrad") rows 2 through 454 all currently hold the
# serial date 45177 (2023-09-08). Update them to 45178 (2023-09-09).
$newDate = (Get-Date -Year 2023 -Month 9 -Day 9 -Hour 0 -Minute 0 -Second 0).Date

$lastRow = 454
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
